$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the default x/y header cells to X/Y (capitalized)
$ws.Range("C1").Value = "X"
$ws.Range("D1").Value = "Y"

# Move selection to D1, matching the saved view state
$ws.Range("D1").Select()
